$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores every value as text in the source data,
# even when it looks like a plain number (e.g. "529.59", "1.01"). Assigning
# such a string straight to .Value would make Excel auto-convert it to a
# Number. Prefixing with a leading apostrophe forces text entry (the
# apostrophe itself is not stored), and resetting .Style to "Normal"
# afterwards keeps the cell formatting identical to the untouched cells
# (values that already contain extra separators, like "57.808.28", are
# not ambiguous and need no special handling).

$ws.Range('D2').Value = '57.808.28'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '3.136.59'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'529.59"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('D6').Value = "'138.61"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.54%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.134.64'
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('E9').Value = '  +2.98%  '
$ws.Range('E10').Value = '  -0.16%  '
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('E12').Value = '  +3.68%  '
$ws.Range('D13').Value = '3.680.21'
$ws.Range('E13').Value = '  +1.20%  '
$ws.Range('E14').Value = '  +2.70%  '
$ws.Range('D15').Value = "'25.56"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.27%  '
$ws.Range('D16').Value = "'0.0000165"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.79%  '
$ws.Range('D17').Value = '57.940.57'
$ws.Range('E17').Value = '  +0.53%  '
$ws.Range('D18').Value = '3.151.19'
$ws.Range('E18').Value = '  +1.59%  '
$ws.Range('E19').Value = '  -2.11%  '
$ws.Range('D20').Value = "'12.73"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('E21').Value = '  -1.13%  '
$ws.Range('D22').Value = "'351.97"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.55%  '
$ws.Range('D23').Value = "'5.80"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.50%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').Value = "'68.54"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.79%  '
$ws.Range('E26').Value = '  -0.63%  '
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('D28').Value = "'1.01"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.85%  '
$ws.Range('D29').Value = '0.0₃0920'
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = "'7.55"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.65%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('D32').Value = "'6.17"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.35%  '
$ws.Range('D33').Value = "'1.88"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.01%  '
$ws.Range('D34').Value = "'21.18"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.23%  '
$ws.Range('E35').Value = '  -0.84%  '
$ws.Range('D36').Value = "'4.99"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.52%  '
$ws.Range('D37').Value = "'157.99"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('D38').Value = "'6.20"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.30%  '
$ws.Range('D39').Value = "'26.48"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.66%  '
$ws.Range('E40').Value = '  -2.31%  '
$ws.Range('D41').Value = "'0.0671"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = "'1.63"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.40%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = "'4.21"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.93%  '
$ws.Range('D44').Value = "'0.705"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.81%  '
$ws.Range('D45').Value = '3.178.06'
$ws.Range('E45').Value = '  +1.12%  '
$ws.Range('E46').Value = '  +4.49%  '
$ws.Range('D47').Value = "'36.65"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.39%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.346.50'
$ws.Range('E48').Value = '  +1.85%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').Value = "'0.999"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('D50').Value = "'0.968"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.23%  '
$ws.Range('D51').Value = "'6.05"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.63%  '
